# This script reproduces a weekly data update: a new observation row is
# inserted into the price table at sheet row 152 (pushing the existing
# rows 152-234 down to 153-235), and the new row is populated with the
# same attributes as the record that used to sit at row 152, except for
# a new date (column D) and a new volume (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 152; everything below (152..234) shifts
# down to (153..235).
$ws.Rows.Item(152).Insert()

# After the insert, row 153 holds what used to be row 152's data. Copy
# those values into the now-empty row 152 so it starts as a duplicate of
# its neighbor, then adjust the two changed fields (Fecha / Volumen).
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(152, $col).Value = $ws.Cells.Item(153, $col).Value()
}

# New date (serial 44529 -> 2021-11-29) and new volume for the inserted row.
$ws.Cells.Item(152, 4).Value = 44529
$ws.Cells.Item(152, 10).Value = 4000
